$wb = $excel.ActiveWorkbook

# --- Contact sheet: add ContactName / ContactType columns (F, G) ---
$contact = $wb.Worksheets.Item("Contact")

$contact.Range("G1").Value = "ContactType"
$contact.Range("G2").Value = "External Contact"
$contact.Range("F1").Value = "ContactName"
$contact.Range("F2").Value = "Test LVContact"

$contact.Range("F1:G1").Font.Bold = $true
$contact.Range("F2").NumberFormat = "@"

$contact.Range("F:F").ColumnWidth = 17.5

$contact.Activate()
$contact.Range("G9").Select() | Out-Null

# --- New sheet: SubscriptionPreferences ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "SubscriptionPreferences"

$newSheet.Range("A1").Value = "Deal Announcements"
$newSheet.Range("A2").Value = "Opt In"
$newSheet.Range("B1").Value = "Event/Conferences"
$newSheet.Range("C1").Value = "General Announcements"
$newSheet.Range("C2").Value = "Opt Out"
$newSheet.Range("D1").Value = "Insights/Content"

$newSheet.Range("B2").Value = "Opt In"
$newSheet.Range("D2").Value = "Opt Out"

$newSheet.Range("A1:D1").Font.Bold = $true

$newSheet.Range("A:A").ColumnWidth = 18.333333333333336
$newSheet.Range("B:B").ColumnWidth = 16.416666666666675
$newSheet.Range("C:C").ColumnWidth = 21.250000000000018

$newSheet.Activate()
$newSheet.Range("G13").Select() | Out-Null
